$p = $ppt.ActivePresentation

# --- 1) Table on slide 5: switch the table to a different built-in table style ---
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shape = $slide5.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{ABCDF956-C958-4ED5-B4DF-287D80FC2B2B}")
    }
}

# --- 2) Re-colour the presentation's theme (Design tab: "Integral" -> "Office Theme") ---
$theme = $p.Slides.Item(1).Master.Theme
$cs = $theme.ThemeColorScheme

$cs.Item(1).RGB  = 0 + 0*256 + 0*65536          # dk1      000000
$cs.Item(2).RGB  = 255 + 255*256 + 255*65536    # lt1      FFFFFF
$cs.Item(3).RGB  = 0x44 + 0x54*256 + 0x6A*65536 # dk2      44546A
$cs.Item(4).RGB  = 0xE7 + 0xE6*256 + 0xE6*65536 # lt2      E7E6E6
$cs.Item(5).RGB  = 0x5B + 0x9B*256 + 0xD5*65536 # accent1  5B9BD5
$cs.Item(6).RGB  = 0xED + 0x7D*256 + 0x31*65536 # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5 + 0xA5*256 + 0xA5*65536 # accent3  A5A5A5
$cs.Item(8).RGB  = 0xFF + 0xC0*256 + 0x00*65536 # accent4  FFC000
$cs.Item(9).RGB  = 0x44 + 0x72*256 + 0xC4*65536 # accent5  4472C4
$cs.Item(10).RGB = 0x70 + 0xAD*256 + 0x47*65536 # accent6  70AD47
$cs.Item(11).RGB = 0x05 + 0x63*256 + 0xC1*65536 # hlink    0563C1
$cs.Item(12).RGB = 0x95 + 0x4F*256 + 0x72*65536 # folHlink 954F72
